$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear old "actual value" / "result" columns (F, G) for all data rows; no longer used
$ws.Range("F2:G14").ClearContents()

# Update existing rows 2-14 and add new rows 15-20 with refreshed test data
$ws.Cells.Item(2,1).Value = 1
$ws.Cells.Item(2,2).Value = 0
$ws.Cells.Item(2,3).Value = 100
$ws.Cells.Item(2,4).Value = 100
$ws.Cells.Item(2,5).Value = "长度超出范围"

$ws.Cells.Item(3,1).Value = 2
$ws.Cells.Item(3,2).Value = 1
$ws.Cells.Item(3,3).Value = 100
$ws.Cells.Item(3,4).Value = 100
$ws.Cells.Item(3,5).Value = "等腰三角形"

$ws.Cells.Item(4,1).Value = 3
$ws.Cells.Item(4,2).Value = 10
$ws.Cells.Item(4,3).Value = 100
$ws.Cells.Item(4,4).Value = 95
$ws.Cells.Item(4,5).Value = "普通三角形"

$ws.Cells.Item(5,1).Value = 4
$ws.Cells.Item(5,2).Value = 100
$ws.Cells.Item(5,3).Value = 100
$ws.Cells.Item(5,4).Value = 100
$ws.Cells.Item(5,5).Value = "等边三角形"

$ws.Cells.Item(6,1).Value = 5
$ws.Cells.Item(6,2).Value = 190
$ws.Cells.Item(6,3).Value = 100
$ws.Cells.Item(6,4).Value = 90
$ws.Cells.Item(6,5).Value = "不构成三角形"

$ws.Cells.Item(7,1).Value = 6
$ws.Cells.Item(7,2).Value = 200
$ws.Cells.Item(7,3).Value = 100
$ws.Cells.Item(7,4).Value = 110
$ws.Cells.Item(7,5).Value = "普通三角形"

$ws.Cells.Item(8,1).Value = 7
$ws.Cells.Item(8,2).Value = 201
$ws.Cells.Item(8,3).Value = 100
$ws.Cells.Item(8,4).Value = 100
$ws.Cells.Item(8,5).Value = "长度超出范围"

$ws.Cells.Item(9,1).Value = 8
$ws.Cells.Item(9,2).Value = 100
$ws.Cells.Item(9,3).Value = 0
$ws.Cells.Item(9,4).Value = 100
$ws.Cells.Item(9,5).Value = "长度超出范围"

$ws.Cells.Item(10,1).Value = 9
$ws.Cells.Item(10,2).Value = 100
$ws.Cells.Item(10,3).Value = 1
$ws.Cells.Item(10,4).Value = 90
$ws.Cells.Item(10,5).Value = "不构成三角形"

$ws.Cells.Item(11,1).Value = 10
$ws.Cells.Item(11,2).Value = 100
$ws.Cells.Item(11,3).Value = 10
$ws.Cells.Item(11,4).Value = 105
$ws.Cells.Item(11,5).Value = "普通三角形"

$ws.Cells.Item(12,1).Value = 11
$ws.Cells.Item(12,2).Value = 100
$ws.Cells.Item(12,3).Value = 190
$ws.Cells.Item(12,4).Value = 110
$ws.Cells.Item(12,5).Value = "普通三角形"

$ws.Cells.Item(13,1).Value = 12
$ws.Cells.Item(13,2).Value = 100
$ws.Cells.Item(13,3).Value = 200
$ws.Cells.Item(13,4).Value = 100
$ws.Cells.Item(13,5).Value = "不构成三角形"

$ws.Cells.Item(14,1).Value = 13
$ws.Cells.Item(14,2).Value = 100
$ws.Cells.Item(14,3).Value = 201
$ws.Cells.Item(14,4).Value = 100
$ws.Cells.Item(14,5).Value = "长度超出范围"

$ws.Cells.Item(15,1).Value = 14
$ws.Cells.Item(15,2).Value = 100
$ws.Cells.Item(15,3).Value = 100
$ws.Cells.Item(15,4).Value = 0
$ws.Cells.Item(15,5).Value = "长度超出范围"

$ws.Cells.Item(16,1).Value = 15
$ws.Cells.Item(16,2).Value = 100
$ws.Cells.Item(16,3).Value = 110
$ws.Cells.Item(16,4).Value = 1
$ws.Cells.Item(16,5).Value = "不构成三角形"

$ws.Cells.Item(17,1).Value = 16
$ws.Cells.Item(17,2).Value = 100
$ws.Cells.Item(17,3).Value = 100
$ws.Cells.Item(17,4).Value = 10
$ws.Cells.Item(17,5).Value = "等腰三角形"

$ws.Cells.Item(18,1).Value = 17
$ws.Cells.Item(18,2).Value = 100
$ws.Cells.Item(18,3).Value = 95
$ws.Cells.Item(18,4).Value = 190
$ws.Cells.Item(18,5).Value = "普通三角形"

$ws.Cells.Item(19,1).Value = 18
$ws.Cells.Item(19,2).Value = 100
$ws.Cells.Item(19,3).Value = 90
$ws.Cells.Item(19,4).Value = 200
$ws.Cells.Item(19,5).Value = "不构成三角形"

$ws.Cells.Item(20,1).Value = 19
$ws.Cells.Item(20,2).Value = 100
$ws.Cells.Item(20,3).Value = 110
$ws.Cells.Item(20,4).Value = 201
$ws.Cells.Item(20,5).Value = "长度超出范围"

# Update selection to reflect the last edited cell location
$ws.Range("F15").Select() | Out-Null
